$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Wow coffee bar"
$ws.Range("C13").Value = "Moreno 2201"
$ws.Range("D13").Value = "[-31.447512460528124, -60.932169283415085]"
$ws.Range("E13").Value = "wowcoffee"
$ws.Range("F13").Value = "wowcoffee2201"

$ws.Range("A27").Font.Underline = $ws.Range("C15").Font.Underline
$ws.Range("B30").Font.Underline = $ws.Range("C15").Font.Underline

$ws.Range("B30").Select()
